$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Deposits")
$ws.Range("K5").Value = "test"
